$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task row: "מספר עובד רנדומלי" (Random employee number) in B9
$ws.Range("B9").Value = "מספר עובד רנדומלי"

# Highlight the two most recent / pending tasks (B8, B9) with a green fill
# (0, 176, 80) == #00B050 packed as an OLE BGR long value
$green = 5287936
$ws.Range("B8:B9").Interior.Color = $green

# Move the active selection to E5 (matches the post-edit workbook state)
$ws.Range("E5").Select() | Out-Null

# Refresh print/page setup so the sheet has an explicit page setup definition
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
